# Encoder ports assigned on X3.
# Applies the edits to the "XMOS Dualchip" worksheet: moves the existing
# "X3 Signal" (col Q) values one column right into col R, and records the
# newly-assigned encoder signal names (ENC0/1/2 _A/_B/_Z) into col P,
# then makes "XMOS Dualchip" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XMOS Dualchip")

# --- Rows where the old column-Q value simply shifts right to column R ---
$ws.Range("R2").Value  = "JP5_9"
$ws.Range("R3").Value  = "JP5_10"
$ws.Range("R12").Value = "JP4_1"
$ws.Range("R13").Value = "JP4_2"
$ws.Range("R27").Value = "JP4_10"
$ws.Range("Q2").ClearContents()
$ws.Range("Q3").ClearContents()
$ws.Range("Q12").ClearContents()
$ws.Range("Q13").ClearContents()
$ws.Range("Q27").ClearContents()

# --- Rows where the old column-Q value shifts to column R, and a new
#     encoder-signal value is recorded in column P ---
$ws.Range("R14").Value = "JP4_3"
$ws.Range("P14").Value = "ENC0_A"
$ws.Range("Q14").ClearContents()

$ws.Range("R15").Value = "JP4_4"
$ws.Range("P15").Value = "ENC1_B"
$ws.Range("Q15").ClearContents()

$ws.Range("R24").Value = "JP4_7"
$ws.Range("P24").Value = "ENC0_Z"
$ws.Range("Q24").ClearContents()

$ws.Range("R25").Value = "JP4_8"
$ws.Range("P25").Value = "ENC1_Z"
$ws.Range("Q25").ClearContents()

$ws.Range("R26").Value = "JP4_9"
$ws.Range("P26").Value = "ENC0_B"
$ws.Range("Q26").ClearContents()

# --- Rows where the old column-P value shifts to column R, and a new
#     encoder-signal value replaces it in column P ---
$ws.Range("R36").Value = "JP5_7"
$ws.Range("P36").Value = "ENC1_A"

$ws.Range("R37").Value = "JP5_8"
$ws.Range("P37").Value = "ENC2_Z"

$ws.Range("R38").Value = "JP5_1"
$ws.Range("P38").Value = "ENC2_A"

$ws.Range("R39").Value = "JP5_2"
$ws.Range("P39").Value = "ENC2_B"

# --- Rows where the old column-P value simply shifts right to column R ---
$ws.Range("R40").Value = "JP5_3"
$ws.Range("P40").ClearContents()

$ws.Range("R41").Value = "JP5_4"
$ws.Range("P41").ClearContents()

# --- View state: "XMOS Dualchip" becomes the active/selected sheet & tab ---
$ws.Activate()
$ws.Range("A7").Select()
$ws.Range("P40").Select()
